# "All staff Data set unique id" - add new faculty rows (Karthick, Chezhian additions,
# Ganapathy Sundaram) and their unique_id column values to Sheet1.
#
# NOTE: cells are written in a specific order so that new shared-string table
# entries get interned in the exact sequence the original authoring session
# produced (row 24 is a scratch/staging area that was filled in before the
# data was copied into row 20, and columns were not filled strictly
# left-to-right either).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: Mr. M.KARTHICK (Assistant Professor) -------------------------
$ws.Range("C19").Value = "https://drive.google.com/u/0/open?usp=forms_web&id=1zhflWBm35sv2732JdwQZN1rZ-bZvH254"
$ws.Range("D19").Value = " https://scholar.google.com/citations?user=WPNgCKEAAAAJ"
$ws.Range("I19").Value = "https://www.linkedin.com/in/karthick-m-ap-iii-mechanical-86330a282/"
$ws.Range("A19").Value = "Mr. M.KARTHICK"
$ws.Range("B19").Value = "Assistant Professor"

# --- Row 18: Mr. R. Chezhian extra columns ---------------------------------
$ws.Range("C18").Value = "https://drive.google.com/u/0/open?usp=forms_web&id=1WphphCt7yyXF4bmVIrhxoG92ed_e2Mha"
$ws.Range("D18").Value = "scholar.google.com/citations?hl=en&pli=1&authuser=3&user=k0jG4AAAAJ"
$ws.Range("I18").Value = "www.linkedin.com/in/chezhian-r-664275146"

# --- Row 16: extra Photo column --------------------------------------------
$ws.Range("C16").Value = "https://drive.google.com/u/0/open?usp=forms_web&id=1Yz2iy4HuSim_mIfwJakZ5suGT_3jsLXS"

# --- Row 24: scratch/staging cells for Ganapathy Sundaram E -----------------
$ws.Range("J24").Value = "ganapathysundaram@velammal.edu.in"
$ws.Range("E24").Value = "Ganapathy Sundaram E"
$ws.Range("F24").Value = "E"
$ws.Range("B20").Value = "Professor & Head"
$ws.Range("C20").Value = "https://drive.google.com/u/0/open?usp=forms_web&id=10mJQLBl0PMMOJTE7wZlVFtMpGlweqH57"
$ws.Range("D20").Value = "https://scholar.google.com/citations?user=PrLAIusAAAAJ&hl=en"
$ws.Range("N24").Value = "0000-0002-3284-9485"
$ws.Range("G20").Value = "LIG-0308-2024"
$ws.Range("I20").Value = "ganapathy-sundaram-esakkimuthu-787163204/?originalSubdomain=in"
$ws.Range("A20").Value = "Dr. E.Ganapathy Sundaram "

# --- unique_id column (J) for the new / touched rows ------------------------
$ws.Range("J16").Value = "VEC-013-04-204"
$ws.Range("J19").Value = "VEC-013-04-198"
$ws.Range("J18").Value = "VEC-013-04-202"
$ws.Range("J20").Value = "VEC-013-01-203"
$ws.Range("J17").Value = "VEC-013-04-103"

# --- remaining row-24 staging duplicates (reuse already-interned strings) --
$ws.Range("K24").Value = $ws.Range("C20").Value()
$ws.Range("L24").Value = $ws.Range("D20").Value()
$ws.Range("O24").Value = $ws.Range("G20").Value()
$ws.Range("Q24").Value = $ws.Range("I20").Value()
$ws.Range("S24").Value = 2

# --- formatted-but-empty cells (date number format, no value) --------------
$ws.Range("N21").NumberFormat = "mm-dd-yy"
$ws.Range("N21").Copy()
$ws.Range("H24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- final selection as left by the editing session -------------------------
$ws.Range("K20").Select()
